# Apply the "Add files via upload" edit: fill in the DFS (H) and Greedy (I)
# execution-time columns for every test row, fix up the BFS (G) value for
# test 5, and move the active selection to I18.
# (Header row: F=A*, G=BFS, H=DFS, I=Greedy, all under "Tempo de Execucao (s)".)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> values for columns G (BFS), H (DFS), I (Greedy)
# ($null means "leave as-is / not changed by this edit")
$rows = @(
    @{ Row = 5;  G = $null;  H = 0.14;  I = 0.135 },
    @{ Row = 6;  G = $null;  H = 0.142; I = 0.14  },
    @{ Row = 7;  G = $null;  H = 0.147; I = 0.151 },
    @{ Row = 8;  G = $null;  H = 0.14;  I = 0.146 },
    @{ Row = 9;  G = 0.159;  H = 0.161; I = 0.171 },
    @{ Row = 10; G = 0.189;  H = 0.182; I = 0.196 },
    @{ Row = 11; G = 0.178;  H = 0.171; I = 0.173 },
    @{ Row = 12; G = 0.14;   H = 0.142; I = 0.159 },
    @{ Row = 13; G = 0.218;  H = 0.222; I = 0.246 },
    @{ Row = 14; G = 0.286;  H = 0.304; I = 0.312 },
    @{ Row = 15; G = 0.235;  H = 0.244; I = 0.268 },
    @{ Row = 16; G = 0.237;  H = 0.242; I = 0.275 },
    @{ Row = 17; G = 0.326;  H = 0.319; I = 0.395 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    if ($null -ne $r.G) {
        $gCell = $ws.Cells.Item($rowNum, 7)   # column G
        $gCell.Value = $r.G
        # Only cells whose value is exactly 0.14 (or that already carried the
        # "0.000" format, like G9) need the explicit number format - otherwise
        # the default General format already renders these values correctly.
        if ($r.G -eq 0.14 -or $rowNum -eq 9) {
            $gCell.NumberFormat = "0.000"
        }
    }

    $hCell = $ws.Cells.Item($rowNum, 8)       # column H (DFS)
    $hCell.Value = $r.H
    if ($r.H -eq 0.14) {
        $hCell.NumberFormat = "0.000"
    }

    $iCell = $ws.Cells.Item($rowNum, 9)       # column I (Greedy)
    $iCell.Value = $r.I
    if ($r.I -eq 0.14) {
        $iCell.NumberFormat = "0.000"
    }
}

# Move the selection like the author left it after editing.
$ws.Range("I18").Select() | Out-Null
